$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = "Labor Office"
$ws.Range("B1").Value = "1-523 - Sanitary Facilities"
$ws.Range("C1").Value = "10/26/2022 1:27:46 PM"
$ws.Range("E1").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/edit.png"
$ws.Range("F1").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Edit/c97f0116-209c-4529-99f5-06b8f7ac890f"
$ws.Range("H1").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/details.png"
$ws.Range("I1").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Details/c97f0116-209c-4529-99f5-06b8f7ac890f"
$ws.Range("K1").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/delete.png"
$ws.Range("L1").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Delete/c97f0116-209c-4529-99f5-06b8f7ac890f"
$ws.Range("M1").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Devices/Details/c03ccdcf-87ef-463f-9a7e-69ff9d5df8f0"
$ws.Range("O1").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/delete.png"
$ws.Range("P1").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Devices/Delete/c03ccdcf-87ef-463f-9a7e-69ff9d5df8f0"

# Row 2
$ws.Range("A2").Value = "Iron Workshop"
$ws.Range("B2").Value = "3-400 - Precast Concrete"
$ws.Range("C2").Value = "10/26/2022 1:27:43 PM"
$ws.Range("E2").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/edit.png"
$ws.Range("F2").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Edit/566d548e-c476-4f24-aa9b-5399d0cba528"
$ws.Range("H2").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/details.png"
$ws.Range("I2").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Details/566d548e-c476-4f24-aa9b-5399d0cba528"
$ws.Range("K2").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/delete.png"
$ws.Range("L2").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Delete/566d548e-c476-4f24-aa9b-5399d0cba528"

# Row 3
$ws.Range("A3").Value = "Terrazzo"
$ws.Range("B3").Value = "13-175 - Ice Rinks"
$ws.Range("C3").Value = "10/26/2022 1:27:32 PM"
$ws.Range("E3").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/edit.png"
$ws.Range("F3").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Edit/86f51092-8bca-46c9-a3fa-7ed6b3d70885"
$ws.Range("H3").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/details.png"
$ws.Range("I3").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Details/86f51092-8bca-46c9-a3fa-7ed6b3d70885"
$ws.Range("K3").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/delete.png"
$ws.Range("L3").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Delete/86f51092-8bca-46c9-a3fa-7ed6b3d70885"

# Row 4
$ws.Range("A4").Value = "Tile Setting Bench"
$ws.Range("B4").Value = "1-570 - Temporary Controls"
$ws.Range("C4").Value = "10/26/2022 1:27:52 PM"
$ws.Range("E4").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/edit.png"
$ws.Range("F4").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Edit/24e10a20-d821-4948-95a1-a8fa4d8413b6"
$ws.Range("H4").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/details.png"
$ws.Range("I4").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Details/24e10a20-d821-4948-95a1-a8fa4d8413b6"
$ws.Range("K4").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/delete.png"
$ws.Range("L4").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Delete/24e10a20-d821-4948-95a1-a8fa4d8413b6"

# Row 5
$ws.Range("A5").Value = "Environmental Office"
$ws.Range("B5").Value = "2-370 - Erosion and Sedimentation Control"
$ws.Range("C5").Value = "10/26/2022 1:27:49 PM"
$ws.Range("E5").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/edit.png"
$ws.Range("F5").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Edit/dea198e1-a276-4e61-bbcb-ad32bc164fde"
$ws.Range("H5").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/details.png"
$ws.Range("I5").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Details/dea198e1-a276-4e61-bbcb-ad32bc164fde"
$ws.Range("K5").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/delete.png"
$ws.Range("L5").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Delete/dea198e1-a276-4e61-bbcb-ad32bc164fde"

# Row 6
$ws.Range("A6").Value = "Boilermaker Room"
$ws.Range("B6").Value = "2-870 - Sculpture/Ornamental"
$ws.Range("C6").Value = "10/26/2022 1:27:35 PM"
$ws.Range("E6").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/edit.png"
$ws.Range("F6").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Edit/771ccd04-4a05-409f-bc98-b4273928ef2a"
$ws.Range("H6").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/details.png"
$ws.Range("I6").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Details/771ccd04-4a05-409f-bc98-b4273928ef2a"
$ws.Range("K6").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/delete.png"
$ws.Range("L6").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Delete/771ccd04-4a05-409f-bc98-b4273928ef2a"

# Row 7
$ws.Range("A7").Value = "Safety Office"
$ws.Range("B7").Value = "17-030 - Bond"
$ws.Range("C7").Value = "10/26/2022 1:27:37 PM"
$ws.Range("E7").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/edit.png"
$ws.Range("F7").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Edit/c9e3b217-595f-4ec7-942c-cbd76eb7aeac"
$ws.Range("H7").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/details.png"
$ws.Range("I7").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Details/c9e3b217-595f-4ec7-942c-cbd76eb7aeac"
$ws.Range("K7").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/delete.png"
$ws.Range("L7").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Delete/c9e3b217-595f-4ec7-942c-cbd76eb7aeac"

# Row 8
$ws.Range("A8").Value = "Stucco Mason Building"
$ws.Range("B8").Value = "2-750 - Concrete Pads and Walks"
$ws.Range("C8").Value = "10/26/2022 1:27:40 PM"
$ws.Range("E8").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/edit.png"
$ws.Range("F8").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Edit/964895a6-cf58-4cb4-b31c-d88f83596eb5"
$ws.Range("H8").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/details.png"
$ws.Range("I8").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Details/964895a6-cf58-4cb4-b31c-d88f83596eb5"
$ws.Range("K8").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/delete.png"
$ws.Range("L8").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Delete/964895a6-cf58-4cb4-b31c-d88f83596eb5"

# Row 9
$ws.Range("A9").Value = "Linemen Assembly"
$ws.Range("B9").Value = "2-825 - Wood Fences and Gates"
$ws.Range("C9").Value = "10/26/2022 1:27:55 PM"
$ws.Range("E9").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/edit.png"
$ws.Range("F9").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Edit/0dd96d1f-9a74-4f13-8c07-e5bc1c2a9231"
$ws.Range("H9").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/details.png"
$ws.Range("I9").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Details/0dd96d1f-9a74-4f13-8c07-e5bc1c2a9231"
$ws.Range("K9").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/delete.png"
$ws.Range("L9").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Delete/0dd96d1f-9a74-4f13-8c07-e5bc1c2a9231"

# Row 10
$ws.Range("A10").Value = "Environmental Office"
$ws.Range("B10").Value = "2-370 - Erosion and Sedimentation Control"
$ws.Range("C10").Value = "10/26/2022 12:47:45 PM"
$ws.Range("E10").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/edit.png"
$ws.Range("F10").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Edit/574a744f-6dc2-46f3-8788-ee7f98e292a2"
$ws.Range("H10").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/details.png"
$ws.Range("I10").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Details/574a744f-6dc2-46f3-8788-ee7f98e292a2"
$ws.Range("K10").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/images/delete.png"
$ws.Range("L10").Value = "https://connectedoffice-devicemanagement.azurewebsites.net/Zones/Delete/574a744f-6dc2-46f3-8788-ee7f98e292a2"

